# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) on Sheet1.
# Row 15 gains a new entrant (Binance-Peg BSC-USD); every coin that used to
# occupy rows 15-51 shifts down one rank, and the former last row (Cronos)
# drops off the bottom of the table. Rows 2-14 keep their coin/link and only
# get refreshed Price / Volume(1h) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 2; D = "67.984.29"; E = "  +1.31%  " },
    @{ Row = 3; D = "2.542.42"; E = "  +0.61%  " },
    @{ Row = 4; E = "  +0.00%  " },
    @{ Row = 5; D = "'592.53"; E = "  +0.47%  " },
    @{ Row = 6; D = "'173.98"; E = "  +0.67%  " },
    @{ Row = 7; E = "  -0.04%  " },
    @{ Row = 8; E = "  -0.31%  " },
    @{ Row = 9; D = "2.541.40"; E = "  +0.59%  " },
    @{ Row = 10; E = "  +0.90%  " },
    @{ Row = 11; E = "  +1.49%  " },
    @{ Row = 12; E = "  -0.90%  " },
    @{ Row = 13; E = "  +0.31%  " },
    @{ Row = 14; D = "'26.57"; E = "  -0.03%  " },
    @{ Row = 15; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "'2.50"; E = "  +150.71%  " },
    @{ Row = 16; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "2.973.78"; E = "  -0.41%  " },
    @{ Row = 17; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "'0.0000178"; E = "  +1.12%  " },
    @{ Row = 18; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "67.864.95"; E = "  +1.24%  " },
    @{ Row = 19; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "2.532.74"; E = "  +0.52%  " },
    @{ Row = 20; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "'11.76"; E = "  +3.53%  " },
    @{ Row = 21; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "'7.97"; E = "  -1.37%  " },
    @{ Row = 22; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "'369.05"; E = "  +4.32%  " },
    @{ Row = 23; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "'4.16"; E = "  -0.52%  " },
    @{ Row = 24; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "'4.60"; E = "  -0.59%  " },
    @{ Row = 25; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "'71.60"; E = "  +2.78%  " },
    @{ Row = 26; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "'1.00"; E = "  +0.03%  " },
    @{ Row = 27; B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "'1.93"; E = "  -3.08%  " },
    @{ Row = 28; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "'10.00"; E = "  +0.51%  " },
    @{ Row = 29; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "2.667.38"; E = "  +0.51%  " },
    @{ Row = 30; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0974"; E = "  -0.32%  " },
    @{ Row = 31; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "'8.47"; E = "  +4.00%  " },
    @{ Row = 32; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "'542.99"; E = "  +2.02%  " },
    @{ Row = 33; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "'1.32"; E = "  -0.29%  " },
    @{ Row = 34; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "'1.87"; E = "  +1.46%  " },
    @{ Row = 35; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "'0.130"; E = "  -1.00%  " },
    @{ Row = 36; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "'0.999"; E = "  -0.09%  " },
    @{ Row = 37; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "'159.40"; E = "  +1.01%  " },
    @{ Row = 38; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "'1.44"; E = "  -1.55%  " },
    @{ Row = 39; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "'19.19"; E = "  +3.02%  " },
    @{ Row = 40; B = "WhiteBITCoin"; C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D = "'18.63"; E = "  +1.00%  " },
    @{ Row = 41; B = "RenderToken"; C = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D = "'5.18"; E = "  +0.90%  " },
    @{ Row = 42; B = "PolygonEcosystemToken"; C = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D = "'0.352"; E = "  -0.56%  " },
    @{ Row = 43; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "'1.79"; E = "  -0.30%  " },
    @{ Row = 44; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "'2.58"; E = "  +3.53%  " },
    @{ Row = 45; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "'1.00"; E = "  +0.24%  " },
    @{ Row = 46; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "'39.28"; E = "  -1.07%  " },
    @{ Row = 47; B = "BabyDogeCoin"; C = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D = "0.0₆0292"; E = "  +4.86%  " },
    @{ Row = 48; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "'148.01"; E = "  -0.72%  " },
    @{ Row = 49; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "'3.72"; E = "  +0.69%  " },
    @{ Row = 50; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "'0.554"; E = "  -0.43%  " },
    @{ Row = 51; B = "Optimism"; C = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"; D = "'1.72"; E = "  +1.59%  " }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in "B", "C", "D", "E") {
        if ($update.ContainsKey($col)) {
            $ws.Range("$col$row").Value = $update[$col]
        }
    }
}

